# Update the "dx" (data element group) column for every data row on the
# "Map" sheet: the old DE_GROUP id (DE_GROUP-QjkuCJf6lCs) is replaced with
# the new mer23 id (DE_GROUP-OuKFZzVk6gr). Row 1 is the header ("dx") and
# is left untouched; data rows run from row 2 through row 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

$ws.Range("B2:B91").Value = "DE_GROUP-OuKFZzVk6gr"

# Mirror the reviewer's final on-screen selection/zoom state from the edit
# session (selecting the column that was just updated).
[void]$ws.Range("B12:B91").Select()
$excel.ActiveWindow.Zoom = 140
